$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph.
$start = $d.Content
$okStart = $start.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx")

# Locate the end of the "(c) 2020 ... Creative Commons Attribution" paragraph.
$end = $d.Content
$okEnd = $end.Find.Execute("Creative Commons Attribution")

if ($okStart -and $okEnd) {
    # Expand to also swallow the blank paragraph right before "Ver no Jupiter ..."
    # (its trailing paragraph mark, one character before $start.Start) and the
    # paragraph mark that ends the copyright paragraph (one character after
    # $end.End), so the whole three paragraphs disappear cleanly.
    $delStart = $start.Start - 1
    $delEnd = $end.End + 1
    $victim = $d.Range($delStart, $delEnd)
    $victim.Delete()
}
